$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so Excel does not
# reinterpret dotted/numeric-looking strings (e.g. "310.30", "1.003")
# as numbers, matching the original inline-string (text) cells.
$touchedCells = @("D2", "E2", "D3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $touchedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.902.43"
$ws.Range("E2").Value = "  -2.18%  "
$ws.Range("D3").Value = "1.833.16"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "310.30"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4610"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("D8").Value = "0.3662"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").Value = "0.07167"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "0.8795"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "0.07833"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "19.62"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.862.00"
$ws.Range("E13").Value = "  -0.05%  "
$ws.Range("D14").Value = "5.337"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "0.000008756"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").Value = "26.930.79"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").Value = "14.49"
$ws.Range("E21").Value = "  -3.01%  "
$ws.Range("D22").Value = "5.014"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Value = "1.979"
$ws.Range("E24").Value = "  +5.65%  "
$ws.Range("D25").Value = "150.69"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("D27").Value = "1.993"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "113.65"
$ws.Range("E28").Value = "  -3.00%  "
$ws.Range("D29").Value = "4.967"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("D30").Value = "0.08839"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "3.130"
$ws.Range("E31").Value = "  +3.58%  "
$ws.Range("D32").Value = "0.7679"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "4.467"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "1.134"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").Value = "2.656"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "1.089"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "0.01937"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").Value = "2.924"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "0.05142"
$ws.Range("E39").Value = "  -2.83%  "
$ws.Range("D40").Value = "6.943"
$ws.Range("E40").Value = "  -3.43%  "
$ws.Range("D41").Value = "0.4977"
$ws.Range("E41").Value = "  -4.66%  "
$ws.Range("D42").Value = "0.1600"
$ws.Range("E42").Value = "  -2.89%  "
$ws.Range("D43").Value = "8.328"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "0.4693"
$ws.Range("E44").Value = "  -3.87%  "
$ws.Range("D45").Value = "10.21"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "102.98"
$ws.Range("E47").Value = "  -0.84%  "
$ws.Range("D48").Value = "1.614"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").Value = "0.06099"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "64.92"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "36.43"
$ws.Range("E51").Value = "  -2.10%  "
